$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - "Conference talk proposal": progress update + status bump
$ws.Range("D5").Value = "mostly done! Just need to write talk description"
$ws.Range("E5").Value = 1

# Row 15 - "Summary slide" progress note
$ws.Range("D15").Value = "lol why"

# Row 16 - "Assessment report" progress note
$ws.Range("D16").Value = "lol fuck"

# Update selection/view to D4 (also clears the old scrolled topLeftCell)
$ws.Range("D4").Select()
